# Update TPM-derived NATMI metrics for the Kng1-Rxfp4 ligand-receptor pair
# sheet. Only the receptor expression (M/N) and the specificity-derived
# columns (O,P,Q,R,S,T) change for rows 2-6; all other data is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 13).Value = 3.442313            # M2
$ws.Cells.Item(2, 14).Value = 10.326939           # N2
$ws.Cells.Item(2, 15).Value = 0.2120737065114005  # O2
$ws.Cells.Item(2, 16).Value = 0.2368526181325179  # P2
$ws.Cells.Item(2, 17).Value = 0.5220290613253333  # Q2
$ws.Cells.Item(2, 18).Value = 4.698261551928      # R2
$ws.Cells.Item(2, 19).Value = 0.2120737065114005  # S2
$ws.Cells.Item(2, 20).Value = 0.2368526181325179  # T2

# Row 3 (only O,P,S,T change; M,N,Q,R stay the same)
$ws.Cells.Item(3, 15).Value = 0.158453673516874   # O3
$ws.Cells.Item(3, 16).Value = 0.1769675649214407  # P3
$ws.Cells.Item(3, 19).Value = 0.158453673516874   # S3
$ws.Cells.Item(3, 20).Value = 0.1769675649214407  # T3

# Row 4
$ws.Cells.Item(4, 13).Value = 2.535264            # M4
$ws.Cells.Item(4, 14).Value = 7.605791999999999   # N4
$ws.Cells.Item(4, 15).Value = 0.1561923141402073  # O4
$ws.Cells.Item(4, 16).Value = 0.174441985971967   # P4
$ws.Cells.Item(4, 17).Value = 0.384474475776      # Q4
$ws.Cells.Item(4, 18).Value = 3.460270281984      # R4
$ws.Cells.Item(4, 19).Value = 0.1561923141402073  # S4
$ws.Cells.Item(4, 20).Value = 0.174441985971967   # T4

# Row 5
$ws.Cells.Item(5, 13).Value = 5.0943505           # M5
$ws.Cells.Item(5, 14).Value = 10.188701           # N5
$ws.Cells.Item(5, 15).Value = 0.3138522826957358  # O5
$ws.Cells.Item(5, 16).Value = 0.2336820724146239  # P5
$ws.Cells.Item(5, 17).Value = 0.7725616495586668  # Q5
$ws.Cells.Item(5, 18).Value = 4.635369897352001   # R5
$ws.Cells.Item(5, 19).Value = 0.3138522826957358  # S5
$ws.Cells.Item(5, 20).Value = 0.2336820724146239  # T5

# Row 6
$ws.Cells.Item(6, 13).Value = 2.587785            # M6
$ws.Cells.Item(6, 14).Value = 7.763355            # N6
$ws.Cells.Item(6, 15).Value = 0.1594280231357824  # O6
$ws.Cells.Item(6, 16).Value = 0.1780557585594505  # P6
$ws.Cells.Item(6, 17).Value = 0.39243932044       # Q6
$ws.Cells.Item(6, 18).Value = 3.53195388396       # R6
$ws.Cells.Item(6, 19).Value = 0.1594280231357824  # S6
$ws.Cells.Item(6, 20).Value = 0.1780557585594505  # T6
